$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

# Row 42
$ws.Range("H42").Value = 300
$ws.Range("I42").Value = 300
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 900
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -670
$ws.Range("N42").ClearContents()

# Row 131
$ws.Range("H131").Value = 2894
$ws.Range("I131").Value = 2894
$ws.Range("K131").Value = 8682
$ws.Range("M131").Value = -3642

# Row 138
$ws.Range("H138").Value = 6349.2285
$ws.Range("J138").Value = 7867.346
$ws.Range("L138").Value = 23602.038
$ws.Range("N138").Value = -33882.038

$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

# Row 102
$ws.Range("H102").Value = 909.2222
$ws.Range("I102").Value = 866.75
$ws.Range("K102").Value = 866.75
$ws.Range("M102").Value = 755.25

# Row 132
$ws.Range("H132").Value = 2567.7
$ws.Range("I132").Value = 1987.1818
$ws.Range("J132").Value = 3277.2222
$ws.Range("K132").Value = 5961.5454
$ws.Range("L132").Value = 9831.6666
$ws.Range("M132").Value = -3431.5454
$ws.Range("N132").Value = -14891.6666

$ws = $wb.Worksheets.Item("BSM")
# Row 76
$ws.Range("H76").Value = 314
$ws.Range("J76").Value = 314
$ws.Range("L76").Value = 314
$ws.Range("N76").Value = -944

# Row 79
$ws.Range("H79").Value = 314
$ws.Range("J79").Value = 314
$ws.Range("L79").Value = 314
$ws.Range("N79").Value = -2498

# Row 134
$ws.Range("H134").Value = 4878.231
$ws.Range("I134").Value = 4909.75
$ws.Range("K134").Value = 14729.25
$ws.Range("M134").Value = -12194.25

$ws = $wb.Worksheets.Item("CRP")
# Row 69
$ws.Range("H69").Value = 9999
$ws.Range("I69").Value = 9999
$ws.Range("K69").Value = 9999
$ws.Range("M69").Value = -9250

# Row 72
$ws.Range("H72").Value = 9999
$ws.Range("I72").Value = 9999
$ws.Range("K72").Value = 29997
$ws.Range("M72").Value = -26253

# Row 132
$ws.Range("H132").Value = 3104.158
$ws.Range("I132").Value = 2099.4546
$ws.Range("J132").Value = 4485.625
$ws.Range("K132").Value = 6298.3638
$ws.Range("L132").Value = 13456.875
$ws.Range("M132").Value = -3768.3638
$ws.Range("N132").Value = -18516.875

# Row 134
$ws.Range("H134").Value = 4423.1577
$ws.Range("I134").Value = 4533.125
$ws.Range("J134").Value = 3836.6667
$ws.Range("K134").Value = 13599.375
$ws.Range("L134").Value = 11510.0001
$ws.Range("M134").Value = -11064.375
$ws.Range("N134").Value = -16580.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 99
$ws.Range("H99").Value = 4833.278
$ws.Range("I99").Value = 1999
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 5997
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = -3751
$ws.Range("N99").Value = -19492

# Row 132
$ws.Range("H132").Value = 7599.6
$ws.Range("J132").Value = 6999.75
$ws.Range("L132").Value = 62997.75
$ws.Range("N132").Value = -68057.75

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 63.5
$ws.Range("I2").Value = 69.666664
$ws.Range("J2").Value = 45
$ws.Range("K2").Value = 69.666664
$ws.Range("L2").Value = 45
$ws.Range("M2").Value = 43.333336
$ws.Range("N2").Value = -271

# Row 122
$ws.Range("H122").Value = 1131.75
$ws.Range("I122").Value = 1131.75
$ws.Range("K122").Value = 3395.25
$ws.Range("M122").Value = -945.25

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1497.25
$ws.Range("J22").Value = 1494
$ws.Range("L22").Value = 1494
$ws.Range("N22").Value = -2084

# Row 27
$ws.Range("H27").Value = 1497.25
$ws.Range("J27").Value = 1494
$ws.Range("L27").Value = 1494
$ws.Range("N27").Value = -1708

# Row 44
$ws.Range("H44").Value = 15000
$ws.Range("J44").Value = 15000
$ws.Range("L44").Value = 15000
$ws.Range("N44").Value = -15912

# Row 82
$ws.Range("H82").Value = 820
$ws.Range("I82").Value = 590
$ws.Range("J82").Value = 896.6667
$ws.Range("K82").Value = 590
$ws.Range("L82").Value = 896.6667
$ws.Range("M82").Value = -229
$ws.Range("N82").Value = -1618.6667

# Row 85
$ws.Range("H85").Value = 820
$ws.Range("I85").Value = 590
$ws.Range("J85").Value = 896.6667
$ws.Range("K85").Value = 590
$ws.Range("L85").Value = 896.6667
$ws.Range("M85").Value = 658
$ws.Range("N85").Value = -3392.6667

# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()

# Row 132
$ws.Range("H132").Value = 2596.6365
$ws.Range("I132").Value = 2229.5
$ws.Range("J132").Value = 3239.125
$ws.Range("K132").Value = 6688.5
$ws.Range("L132").Value = 9717.375
$ws.Range("M132").Value = -4158.5
$ws.Range("N132").Value = -14777.375

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 7672.3335
$ws.Range("I62").Value = 4711.6665
$ws.Range("J62").Value = 10633
$ws.Range("K62").Value = 4711.6665
$ws.Range("L62").Value = 10633
$ws.Range("M62").Value = -4087.6665
$ws.Range("N62").Value = -11881

# Row 65
$ws.Range("H65").Value = 7672.3335
$ws.Range("I65").Value = 4711.6665
$ws.Range("J65").Value = 10633
$ws.Range("K65").Value = 23558.3325
$ws.Range("L65").Value = 53165
$ws.Range("M65").Value = -20438.3325
$ws.Range("N65").Value = -59405

# Row 81
$ws.Range("H81").Value = 20494.75
$ws.Range("J81").Value = 22002
$ws.Range("L81").Value = 44004
$ws.Range("N81").Value = -46126

# Row 84
$ws.Range("H84").Value = 20494.75
$ws.Range("J84").Value = 22002
$ws.Range("L84").Value = 220020
$ws.Range("N84").Value = -230628

# Row 132
$ws.Range("H132").Value = 2835.95
$ws.Range("I132").Value = 2069.75
$ws.Range("J132").Value = 3346.75
$ws.Range("K132").Value = 6209.25
$ws.Range("L132").Value = 10040.25
$ws.Range("M132").Value = -3679.25
$ws.Range("N132").Value = -15100.25
